# Applies the data reshuffle described in the commit diff to the
# "dataDefinition" worksheet. The block of rows describing the
# "VerifiedDiagnosis" diagnosis (originally rows 14-15) is moved to the
# top of the data table (rows 2-3), and the remaining Observation blocks
# are rearranged among rows 4-15 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("dataDefinition")

# Row 2 (was CardiacArrestAfter12h / Observation) -> VerifiedDiagnosis / Diagnosis
$ws.Range("A2").Value = "Diagnosis"
$ws.Range("B2").Value = "VerifiedDiagnosis"
# Leading apostrophe forces text storage (matches original "True" text cells
# like H18/H19) instead of Excel auto-coercing the literal to a boolean.
$ws.Range("H2").Value = "'True"

# Row 3 (was ResultValue for CardiacArrestAfter12h) -> DiagnosisType
$ws.Range("D3").Value = "DiagnosisType"
$ws.Range("F3").Value = "Den endelige diagnose verificeret af kardiolog ved udskrivelse."
$ws.Range("G3").Value = 'Enums/Udfald: | "di200" | "di213a" | "di213b" | "di213c" | "di214" | "ikke_aks" | '
$ws.Range("H3").Value = "'True"

# Row 4 (was LVEF) -> CardiogenicShockAfter12h
$ws.Range("B4").Value = "CardiogenicShockAfter12h"

# Row 5 (was ResultValue for LVEF) -> ResultValue for CardiogenicShockAfter12h
$ws.Range("E5").Value = "str, Enum"
$ws.Range("F5").Value = "Kardiogent shock senere end 12 timer efter indlæggelse."
$ws.Range("G5").Value = 'Enums/Udfald: | "ja" | "nej" | "ikke_relevant" | '

# Row 6 (was CardiacEcho) -> CardiacArrestWithin12h
$ws.Range("B6").Value = "CardiacArrestWithin12h"

# Row 7 (was StatusCode for CardiacEcho) -> ResultValue for CardiacArrestWithin12h
$ws.Range("D7").Value = "ResultValue"
$ws.Range("F7").Value = "Se webservice dokumentation."

# Row 8 (was CardiogenicShockAfter12h) -> AcuteHeartFailureKillipClass
$ws.Range("B8").Value = "AcuteHeartFailureKillipClass"

# Row 9 (was ResultValue for CardiogenicShockAfter12h) -> ResultValue for AcuteHeartFailureKillipClass
$ws.Range("F9").Value = "Killip klasse indenfor 12 timer efter indlæggelse."
$ws.Range("G9").Value = 'Enums/Udfald: | "killip_klasse1" | "killip_klasse2" | "killip_klasse3" | "killip_klasse4" | '

# Row 10 (was AcuteHeartFailureKillipClass) -> LVEF
$ws.Range("B10").Value = "LVEF"

# Row 11 (was ResultValue for AcuteHeartFailureKillipClass) -> ResultValue for LVEF
$ws.Range("E11").Value = "Integer"
$ws.Range("F11").Value = "Se webservice dokumentation."
$ws.Range("G11").Value = "Greater than or equal to: 0 | Less than or equal to: 100"

# Row 12 (was CardiacArrestWithin12h) -> CardiacArrestAfter12h
$ws.Range("B12").Value = "CardiacArrestAfter12h"

# Row 13 (was ResultValue for CardiacArrestWithin12h) -> ResultValue for CardiacArrestAfter12h
$ws.Range("F13").Value = "Hjertestop senere end 12 timer efter indlæggelse."
$ws.Range("G13").Value = 'Enums/Udfald: | "ja" | "nej" | "ikke_relevant" | '

# Row 14 (was Diagnosis / VerifiedDiagnosis) -> Observation / CardiacEcho
$ws.Range("A14").Value = "Observation"
$ws.Range("B14").Value = "CardiacEcho"
$ws.Range("H14").ClearContents()

# Row 15 (was DiagnosisType) -> StatusCode for CardiacEcho
$ws.Range("D15").Value = "StatusCode"
$ws.Range("F15").Value = "Er Ekkokardiografi udført?"
$ws.Range("G15").Value = 'Enums/Udfald: | "ja" | "nej" | '
$ws.Range("H15").ClearContents()
